# alt-chunk.docx fixture update:
#  - outer paragraph 1 gets new text and now ends with a (continuous)
#    section break of its own (the "inner" section), carrying the page
#    size that LO previously dropped for inner altChunk sections.
#  - outer paragraph 2 gets new text.
#  - the trailing (outer/body) sectPr becomes a *continuous* section with
#    an updated page size/margins/columns to match the inner one, so no
#    spurious page break appears before the altChunk content.

$d = $word.ActiveDocument

# 1) First paragraph: new text, plus a new section break (w:pPr/w:sectPr)
#    attached to its paragraph mark.
$p1 = $d.Paragraphs(1)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:pPr>' +
                '<w:sectPr>' +
                  '<w:pgSz w:w="11906" w:h="16838"/>' +
                  '<w:pgMar w:top="1417" w:right="1417" w:bottom="1417" w:left="1417" w:header="708" w:footer="708" w:gutter="0"/>' +
                  '<w:pgNumType w:start="1"/>' +
                  '<w:cols w:space="708"/>' +
                  '<w:docGrid w:linePitch="360"/>' +
                '</w:sectPr>' +
              '</w:pPr>' +
              '<w:r><w:t>outer, before sect break</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$p1.Range.InsertXML($xml1)

# 2) Second paragraph: just the new text.
$p2 = $d.Paragraphs(2)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p><w:r><w:t>outer, after sect break</w:t></w:r></w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'
$p2.Range.InsertXML($xml2)

# 3) Trailing/body sectPr (after the altChunk): make it a continuous
#    section and give it the same page size/margins/columns as above.
# wdSectionContinuous = 0 -> serializes as <w:type w:val="continuous"/>
$last = $d.Sections($d.Sections.Count)
$ps = $last.PageSetup
$ps.SectionStart = 0
# points = twips / 20: 11906/20, 16838/20, 708/20, 708/20, 708/20
$ps.PageWidth = 595.3
$ps.PageHeight = 841.9
$ps.HeaderDistance = 35.4
$ps.FooterDistance = 35.4
$ps.TextColumns.Spacing = 35.4
